$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 17: date, task (B), and result (D) - matches the daily-routine log pattern
# Copy the date cell's style from the row above (A16) so the same date
# number format (s="1") is reused instead of minting a new numFmt.
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17").Value = 43838

$ws.Range("B17").Value = "design the front page"
$ws.Range("D17").Value = "Design the front page"
